$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3980.8823
$ws.Range("I19").Value = 7457.2856
$ws.Range("J19").Value = 1547.4
$ws.Range("K19").Value = 7457.2856
$ws.Range("L19").Value = 1547.4
$ws.Range("M19").Value = -7282.2856
$ws.Range("N19").Value = -1897.4
$ws.Range("H41").Value = 3626.7
$ws.Range("J41").Value = 8571.333000000001
$ws.Range("L41").Value = 8571.333000000001
$ws.Range("N41").Value = -9451.333000000001
$ws.Range("H64").Value = 4045.805
$ws.Range("I64").Value = 3706.4194
$ws.Range("J64").Value = 5097.9
$ws.Range("K64").Value = 3706.4194
$ws.Range("L64").Value = 5097.9
$ws.Range("M64").Value = -3458.4194
$ws.Range("N64").Value = -5593.9
$ws.Range("H67").Value = 4045.805
$ws.Range("I67").Value = 3706.4194
$ws.Range("J67").Value = 5097.9
$ws.Range("K67").Value = 3706.4194
$ws.Range("L67").Value = 5097.9
$ws.Range("M67").Value = -2848.4194
$ws.Range("N67").Value = -6813.9
$ws.Range("H74").Value = 4822.963
$ws.Range("I74").Value = 4786.154
$ws.Range("J74").Value = 4857.143
$ws.Range("K74").Value = 4786.154
$ws.Range("L74").Value = 4857.143
$ws.Range("M74").Value = -3850.154
$ws.Range("N74").Value = -6729.143
$ws.Range("H76").Value = 22735880
$ws.Range("I76").Value = 41678092
$ws.Range("J76").Value = 5226.1
$ws.Range("K76").Value = 41678092
$ws.Range("L76").Value = 5226.1
$ws.Range("M76").Value = -41677777
$ws.Range("N76").Value = -5856.1
$ws.Range("H77").Value = 4822.963
$ws.Range("I77").Value = 4786.154
$ws.Range("J77").Value = 4857.143
$ws.Range("K77").Value = 23930.77
$ws.Range("L77").Value = 24285.715
$ws.Range("M77").Value = -19250.77
$ws.Range("N77").Value = -33645.715
$ws.Range("H79").Value = 22735880
$ws.Range("I79").Value = 41678092
$ws.Range("J79").Value = 5226.1
$ws.Range("K79").Value = 41678092
$ws.Range("L79").Value = 5226.1
$ws.Range("M79").Value = -41677000
$ws.Range("N79").Value = -7410.1
$ws.Range("H116").Value = 73784.734
$ws.Range("I116").Value = 99253.73
$ws.Range("J116").Value = 3745
$ws.Range("K116").Value = 99253.73
$ws.Range("L116").Value = 3745
$ws.Range("M116").Value = -95811.73
$ws.Range("N116").Value = -10629
$ws.Range("H125").Value = 1870.4706
$ws.Range("I125").Value = 943.6667
$ws.Range("J125").Value = 2376
$ws.Range("K125").Value = 8493.0003
$ws.Range("L125").Value = 21384
$ws.Range("M125").Value = -6033.0003
$ws.Range("N125").Value = -26304
$ws.Range("H137").Value = 2965.7
$ws.Range("I137").Value = 3007.348
$ws.Range("J137").Value = 2828.8572
$ws.Range("K137").Value = 9022.044
$ws.Range("L137").Value = 8486.571599999999
$ws.Range("M137").Value = -6472.044
$ws.Range("N137").Value = -13586.5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1285700
$ws.Range("I32").Value = 1518844
$ws.Range("K32").Value = 1518844
$ws.Range("M32").Value = -1518557
$ws.Range("H63").Value = 2657.25
$ws.Range("J63").Value = 1972.5
$ws.Range("L63").Value = 1972.5
$ws.Range("N63").Value = -3344.5
$ws.Range("H66").Value = 2657.25
$ws.Range("J66").Value = 1972.5
$ws.Range("L66").Value = 9862.5
$ws.Range("N66").Value = -16726.5
$ws.Range("H97").Value = 924.6774
$ws.Range("I97").Value = 1043.9584
$ws.Range("J97").Value = 515.7143
$ws.Range("K97").Value = 1043.9584
$ws.Range("L97").Value = 515.7143
$ws.Range("M97").Value = -547.9584
$ws.Range("N97").Value = -1507.7143
$ws.Range("H110").Value = 3097.5356
$ws.Range("I110").Value = 3119.6667
$ws.Range("J110").Value = 2500
$ws.Range("K110").Value = 3119.6667
$ws.Range("L110").Value = 2500
$ws.Range("M110").Value = -1074.6667
$ws.Range("N110").Value = -6590
$ws.Range("H132").Value = 28016.902
$ws.Range("I132").Value = 59612.11
$ws.Range("J132").Value = 3290.2173
$ws.Range("K132").Value = 178836.33
$ws.Range("L132").Value = 9870.651899999999
$ws.Range("M132").Value = -176306.33
$ws.Range("N132").Value = -14930.6519

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3546.5557
$ws.Range("I86").Value = 5014.0713
$ws.Range("J86").Value = 1966.1538
$ws.Range("K86").Value = 5014.0713
$ws.Range("L86").Value = 1966.1538
$ws.Range("M86").Value = -3891.0713
$ws.Range("N86").Value = -4212.1538
$ws.Range("H89").Value = 3546.5557
$ws.Range("I89").Value = 5014.0713
$ws.Range("J89").Value = 1966.1538
$ws.Range("K89").Value = 25070.3565
$ws.Range("L89").Value = 9830.769
$ws.Range("M89").Value = -19454.3565
$ws.Range("N89").Value = -21062.769

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2497.425
$ws.Range("I31").Value = 1093.4517
$ws.Range("J31").Value = 7333.3335
$ws.Range("K31").Value = 1093.4517
$ws.Range("L31").Value = 7333.3335
$ws.Range("M31").Value = -798.4517000000001
$ws.Range("N31").Value = -7923.3335
$ws.Range("H34").Value = 2497.425
$ws.Range("I34").Value = 1093.4517
$ws.Range("J34").Value = 7333.3335
$ws.Range("K34").Value = 1093.4517
$ws.Range("L34").Value = 7333.3335
$ws.Range("M34").Value = -891.4517000000001
$ws.Range("N34").Value = -7737.3335
$ws.Range("H58").Value = 3586.4375
$ws.Range("I58").Value = 4712.84
$ws.Range("J58").Value = 2362.087
$ws.Range("K58").Value = 4712.84
$ws.Range("L58").Value = 2362.087
$ws.Range("M58").Value = -4509.84
$ws.Range("N58").Value = -2768.087
$ws.Range("H62").Value = 2781026
$ws.Range("I62").Value = 6946686
$ws.Range("J62").Value = 3919.4167
$ws.Range("K62").Value = 6946686
$ws.Range("L62").Value = 3919.4167
$ws.Range("M62").Value = -6946062
$ws.Range("N62").Value = -5167.4167
$ws.Range("H65").Value = 2781026
$ws.Range("I65").Value = 6946686
$ws.Range("J65").Value = 3919.4167
$ws.Range("K65").Value = 34733430
$ws.Range("L65").Value = 19597.0835
$ws.Range("M65").Value = -34730310
$ws.Range("N65").Value = -25837.0835
$ws.Range("H94").Value = 5408.15
$ws.Range("I94").Value = 1384.3334
$ws.Range("K94").Value = 1384.3334
$ws.Range("M94").Value = -933.3334
$ws.Range("H136").Value = 3586.4375
$ws.Range("I136").Value = 4712.84
$ws.Range("J136").Value = 2362.087
$ws.Range("K136").Value = 14138.52
$ws.Range("L136").Value = 7086.261
$ws.Range("M136").Value = -11588.52
$ws.Range("N136").Value = -12186.261

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 3999.2
$ws.Range("J93").Value = 4499
$ws.Range("L93").Value = 13497
$ws.Range("N93").Value = -17241
$ws.Range("H114").Value = 1845.7142
$ws.Range("I114").Value = 2622.3333
$ws.Range("J114").Value = 810.2222
$ws.Range("K114").Value = 7866.999899999999
$ws.Range("L114").Value = 2430.6666
$ws.Range("M114").Value = -4612.999899999999
$ws.Range("N114").Value = -8938.6666
$ws.Range("H129").Value = 3212.9443
$ws.Range("I129").Value = 3603
$ws.Range("J129").Value = 2725.375
$ws.Range("K129").Value = 10809
$ws.Range("L129").Value = 8176.125
$ws.Range("M129").Value = -5809
$ws.Range("N129").Value = -18176.125
$ws.Range("H131").Value = 2187.5588
$ws.Range("J131").Value = 1816.7931
$ws.Range("L131").Value = 5450.379300000001
$ws.Range("N131").Value = -15530.3793

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7375.561
$ws.Range("I70").Value = 4504.6665
$ws.Range("J70").Value = 10390
$ws.Range("K70").Value = 4504.6665
$ws.Range("L70").Value = 10390
$ws.Range("M70").Value = -4234.6665
$ws.Range("N70").Value = -10930
$ws.Range("H73").Value = 7375.561
$ws.Range("I73").Value = 4504.6665
$ws.Range("J73").Value = 10390
$ws.Range("K73").Value = 4504.6665
$ws.Range("L73").Value = 10390
$ws.Range("M73").Value = -3568.6665
$ws.Range("N73").Value = -12262
$ws.Range("H80").Value = 4661.25
$ws.Range("I80").Value = 5150.1787
$ws.Range("K80").Value = 5150.1787
$ws.Range("M80").Value = -4152.1787
$ws.Range("H83").Value = 4661.25
$ws.Range("I83").Value = 5150.1787
$ws.Range("K83").Value = 25750.8935
$ws.Range("M83").Value = -20758.8935
$ws.Range("H113").Value = 1945
$ws.Range("I113").Value = 1747.4706
$ws.Range("J113").Value = 2504.6667
$ws.Range("K113").Value = 1747.4706
$ws.Range("L113").Value = 2504.6667
$ws.Range("M113").Value = 422.5293999999999
$ws.Range("N113").Value = -6844.6667
$ws.Range("H122").Value = 1595.1666
$ws.Range("I122").Value = 1168.1111
$ws.Range("J122").Value = 2022.2222
$ws.Range("K122").Value = 3504.3333
$ws.Range("L122").Value = 6066.6666
$ws.Range("M122").Value = -1054.3333
$ws.Range("N122").Value = -10966.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2965.2727
$ws.Range("I7").Value = 2771.75
$ws.Range("J7").Value = 3481.3333
$ws.Range("K7").Value = 2771.75
$ws.Range("L7").Value = 3481.3333
$ws.Range("M7").Value = -2659.75
$ws.Range("N7").Value = -3705.3333
$ws.Range("H55").Value = 232.8125
$ws.Range("I55").Value = 168.88889
$ws.Range("K55").Value = 168.88889
$ws.Range("M55").Value = 4.111109999999996
$ws.Range("H122").Value = 3026.9412
$ws.Range("I122").Value = 2662
$ws.Range("J122").Value = 3437.5
$ws.Range("K122").Value = 7986
$ws.Range("L122").Value = 10312.5
$ws.Range("M122").Value = -5536
$ws.Range("N122").Value = -15212.5
$ws.Range("H126").Value = 2965.2727
$ws.Range("I126").Value = 2771.75
$ws.Range("J126").Value = 3481.3333
$ws.Range("K126").Value = 8315.25
$ws.Range("L126").Value = 10443.9999
$ws.Range("M126").Value = -5845.25
$ws.Range("N126").Value = -15383.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 40004370
$ws.Range("I132").Value = 62504384
$ws.Range("J132").Value = 4333
$ws.Range("K132").Value = 187513152
$ws.Range("L132").Value = 12999
$ws.Range("M132").Value = -187510622
$ws.Range("N132").Value = -18059
